$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as text so values like "41.926.38" or "1.00"
# are preserved exactly as strings rather than being reinterpreted as numbers/dates.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "42.032.92"
$ws.Range("E2").Value = "  -1.80%  "

# Row 3
$ws.Range("D3").Value = "2.248.16"
$ws.Range("E3").Value = "  -3.70%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "297.49"
$ws.Range("E5").Value = "  -2.95%  "

# Row 6
$ws.Range("D6").Value = "94.17"
$ws.Range("E6").Value = "  -6.38%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.495"
$ws.Range("E7").Value = "  -2.61%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -4.41%  "

# Row 10
$ws.Range("D10").Value = "32.58"
$ws.Range("E10").Value = "  -6.90%  "

# Row 11
$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").Value = "  -2.23%  "

# Row 12
$ws.Range("D12").Value = "48.66"
$ws.Range("E12").Value = "  -6.64%  "

# Row 13
$ws.Range("E13").Value = "  +0.36%  "

# Row 14
$ws.Range("D14").Value = "6.62"
$ws.Range("E14").Value = "  -2.82%  "

# Row 15
$ws.Range("D15").Value = "2.596.65"
$ws.Range("E15").Value = "  -3.85%  "

# Row 16
$ws.Range("D16").Value = "15.27"
$ws.Range("E16").Value = "  -2.24%  "

# Row 17
$ws.Range("D17").Value = "2.245.12"
$ws.Range("E17").Value = "  -2.74%  "

# Row 18
$ws.Range("D18").Value = "0.770"
$ws.Range("E18").Value = "  -3.41%  "

# Row 19
$ws.Range("D19").Value = "41.936.94"
$ws.Range("E19").Value = "  -1.87%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0886"
$ws.Range("E20").Value = "  -2.13%  "

# Row 21
$ws.Range("D21").Value = "11.34"
$ws.Range("E21").Value = "  -2.51%  "

# Row 22
$ws.Range("D22").Value = "5.95"
$ws.Range("E22").Value = "  -5.32%  "

# Row 23
$ws.Range("D23").Value = "66.05"
$ws.Range("E23").Value = "  -2.45%  "

# Row 24
$ws.Range("D24").Value = "231.73"
$ws.Range("E24").Value = "  -1.89%  "

# Row 25
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").Value = "  -4.53%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("E27").Value = "  -4.63%  "

# Row 28
$ws.Range("D28").Value = "23.72"
$ws.Range("E28").Value = "  -4.72%  "

# Row 29
$ws.Range("E29").Value = "  -6.87%  "

# Row 30
$ws.Range("D30").Value = "165.34"
$ws.Range("E30").Value = "  +3.61%  "

# Row 31
$ws.Range("D31").Value = "33.46"
$ws.Range("E31").Value = "  -4.39%  "

# Row 32
$ws.Range("D32").Value = "9.02"
$ws.Range("E32").Value = "  -3.58%  "

# Row 33
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
$ws.Range("D34").Value = "4.91"
$ws.Range("E34").Value = "  -4.14%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0690"
$ws.Range("E35").Value = "  -5.23%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.32"
$ws.Range("E36").Value = "  -5.23%  "

# Row 37
$ws.Range("D37").Value = "4.31"
$ws.Range("E37").Value = "  -6.03%  "

# Row 38
$ws.Range("D38").Value = "2.78"
$ws.Range("E38").Value = "  -5.86%  "

# Row 39
$ws.Range("D39").Value = "15.80"
$ws.Range("E39").Value = "  -9.25%  "

# Row 40
$ws.Range("D40").Value = "0.0981"
$ws.Range("E40").Value = "  -4.81%  "

# Row 41
$ws.Range("E41").Value = "  -2.88%  "

# Row 42
$ws.Range("D42").Value = "1.70"
$ws.Range("E42").Value = "  -8.51%  "

# Row 43
$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  +3.01%  "

# Row 44
$ws.Range("D44").Value = "1.936.07"
$ws.Range("E44").Value = "  -3.91%  "

# Row 45
$ws.Range("E45").Value = "  -3.07%  "

# Row 46
$ws.Range("D46").Value = "17.30"
$ws.Range("E46").Value = "  -8.46%  "

# Row 47
$ws.Range("E47").Value = "  -8.49%  "

# Row 48
$ws.Range("E48").Value = "  -6.14%  "

# Row 49
$ws.Range("D49").Value = "2.82"
$ws.Range("E49").Value = "  -3.14%  "

# Row 50
$ws.Range("D50").Value = "2.476.61"
$ws.Range("E50").Value = "  -3.08%  "

# Row 51
$ws.Range("D51").Value = "51.95"
$ws.Range("E51").Value = "  -7.17%  "
